# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (column E) and its corresponding "Fecha" (column F)
# values on rows 16 and 18 were swapped: row 16 now carries period 1802
# (date 09/06/2016) and row 18 now carries period 1809 (date 08/08/2017).
# Row 17 (period 1808) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap "Periodo Mora" text values between rows 16 and 18
$ws.Range("E16").Value = "1802"
$ws.Range("E18").Value = "1809"

# Swap the matching date ("Fecha") values between rows 16 and 18
$ws.Range("F16").Value = 42530
$ws.Range("F18").Value = 42955
